$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("BO4").Value = "CP_RECEPTOR"
$ws.Range("BP4").Value = "HORARIO"
